$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The document currently has, in order:
#   ...
#   3. "Tag the last commit to master with this release number (e.g., Release0.1.2)."
#   4. "Change the define in CLDefinitions from cliff servers to staging servers"
#   5. "Change properties in CloudApiPublic, ..."
#   ...
#
# The target order is:
#   ...
#   3. "Change the define in CLDefinitions from cliff servers to staging servers"
#   4. "Test the application under Visual Studio to make sure it is working."
#      (carries the _GoBack bookmark, which used to sit near the very end
#       of the document)
#   5. "If anything needed to be changed, commit and push to master again."
#   6. "Tag the last commit to master with this release number (e.g., Release0.1.2)."
#   7. "Change properties in CloudApiPublic, ..."
#   ...
# ------------------------------------------------------------------

# Paragraph holding "Change the define in CLDefinitions ..." (currently #4).
$changeDefinePara = $d.Paragraphs(4).Range

# Create three fresh, empty paragraphs right after it.
$changeDefinePara.InsertParagraphAfter()
$d.Paragraphs(5).Range.InsertParagraphAfter()
$d.Paragraphs(6).Range.InsertParagraphAfter()

# --- Paragraph 5: "Test the application under Visual Studio to make sure it is working." ---
$p5Start = $d.Paragraphs(5).Range.Start
$p5FirstRun = $d.Range($p5Start, $p5Start)
$p5FirstRun.InsertAfter("Test the application under Visual Studio ")

$p5MidPos = $d.Paragraphs(5).Range.End - 1
$p5SecondRun = $d.Range($p5MidPos, $p5MidPos)
$p5SecondRun.InsertAfter("to make sure it is working.")

# --- Paragraph 6: "If anything needed to be changed, commit and push to master again." ---
$d.Paragraphs(6).Range.Text = "If anything needed to be changed, commit and push to master again."

# --- Paragraph 7: "Tag the last commit to master with this release number (e.g., Release0.1.2)." ---
$d.Paragraphs(7).Range.Text = "Tag the last commit to master with this release number (e.g., Release0.1.2)."

# Move the "_GoBack" bookmark onto the start of the new "Test the application..."
# paragraph (this also removes it from its old location near the end of the
# document, since a bookmark name is unique within the document).
$goBackPos = $d.Paragraphs(5).Range.Start
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Finally, delete the original "Tag the last commit..." paragraph (#3), which
# has now been superseded by the copy re-inserted as paragraph 7 above.
$d.Paragraphs(3).Range.Delete()
